$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column (D) retains its original text formatting so that
# numeric-looking strings (e.g. "226.50") are not coerced into numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value2 = "34.166.29"
$ws.Range("E2").Value2 = "  +0.67%  "
$ws.Range("D3").Value2 = "1.788.17"
$ws.Range("E3").Value2 = "  +0.63%  "
$ws.Range("E4").Value2 = "  +0.08%  "
$ws.Range("D5").Value2 = "226.50"
$ws.Range("E5").Value2 = "  +0.54%  "
$ws.Range("E6").Value2 = "  -1.00%  "
$ws.Range("E7").Value2 = "  +0.08%  "
$ws.Range("D8").Value2 = "31.84"
$ws.Range("E8").Value2 = "  -0.63%  "
$ws.Range("E9").Value2 = "  +1.07%  "
$ws.Range("E10").Value2 = "  -1.29%  "
$ws.Range("E11").Value2 = "  +0.71%  "
$ws.Range("D12").Value2 = "2.046.54"
$ws.Range("E12").Value2 = "  +0.60%  "
$ws.Range("B13").Value2 = "Chainlink"
$ws.Range("C13").Value2 = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D13").Value2 = "11.05"
$ws.Range("E13").Value2 = "  +0.72%  "
$ws.Range("B14").Value2 = "WrappedEther"
$ws.Range("C14").Value2 = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D14").Value2 = "1.796.69"
$ws.Range("E14").Value2 = "  +1.37%  "
$ws.Range("D15").Value2 = "34.131.61"
$ws.Range("E15").Value2 = "  +0.59%  "
$ws.Range("E16").Value2 = "  +0.74%  "
$ws.Range("E17").Value2 = "  +0.53%  "
$ws.Range("D18").Value2 = "68.24"
$ws.Range("E18").Value2 = "  +1.15%  "
$ws.Range("D19").Value2 = "247.21"
$ws.Range("E19").Value2 = "  +2.71%  "
$ws.Range("E20").Value2 = "  -0.38%  "
$ws.Range("D22").Value2 = "10.84"
$ws.Range("E22").Value2 = "  +1.64%  "
$ws.Range("E23").Value2 = "  +0.44%  "
$ws.Range("D24").Value2 = "2.04"
$ws.Range("E24").Value2 = "  -0.42%  "
$ws.Range("D25").Value2 = "161.01"
$ws.Range("E25").Value2 = "  +0.85%  "
$ws.Range("E26").Value2 = "  +1.36%  "
$ws.Range("D27").Value2 = "16.33"
$ws.Range("E27").Value2 = "  +0.69%  "
$ws.Range("E28").Value2 = "  +0.82%  "
$ws.Range("E29").Value2 = "  +0.11%  "
$ws.Range("E30").Value2 = "  +0.10%  "
$ws.Range("D31").Value2 = "0.0518"
$ws.Range("E31").Value2 = "  +1.29%  "
$ws.Range("D32").Value2 = "3.66"
$ws.Range("E32").Value2 = "  +0.98%  "
$ws.Range("E33").Value2 = "  +2.53%  "
$ws.Range("E34").Value2 = "  -0.23%  "
$ws.Range("D35").Value2 = "1.447.23"
$ws.Range("E35").Value2 = "  +4.15%  "
$ws.Range("B36").Value2 = "RenderToken"
$ws.Range("C36").Value2 = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D36").Value2 = "2.43"
$ws.Range("E36").Value2 = "  +8.50%  "
$ws.Range("B37").Value2 = "ImmutableX"
$ws.Range("C37").Value2 = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D37").Value2 = "0.648"
$ws.Range("E37").Value2 = "  -1.46%  "
$ws.Range("D38").Value2 = "0.0191"
$ws.Range("E38").Value2 = "  +2.91%  "
$ws.Range("E39").Value2 = "  +0.06%  "
$ws.Range("B40").Value2 = "HuobiToken"
$ws.Range("C40").Value2 = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D40").Value2 = "2.38"
$ws.Range("E40").Value2 = "  +1.04%  "
$ws.Range("B41").Value2 = "Aave"
$ws.Range("C41").Value2 = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D41").Value2 = "80.36"
$ws.Range("E41").Value2 = "  +3.40%  "
$ws.Range("D42").Value2 = "0.923"
$ws.Range("E42").Value2 = "  +1.48%  "
$ws.Range("E43").Value2 = "  +1.28%  "
$ws.Range("E44").Value2 = "  +0.72%  "
$ws.Range("B45").Value2 = "FraxShare"
$ws.Range("C45").Value2 = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D45").Value2 = "6.05"
$ws.Range("E45").Value2 = "  +3.32%  "
$ws.Range("B46").Value2 = "Kaspa"
$ws.Range("C46").Value2 = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D46").Value2 = "0.0508"
$ws.Range("E46").Value2 = "  +2.27%  "
$ws.Range("E47").Value2 = "  -0.21%  "
$ws.Range("E48").Value2 = "  -7.65%  "
$ws.Range("D49").Value2 = "1.947.96"
$ws.Range("E49").Value2 = "  +0.85%  "
$ws.Range("D50").Value2 = "105.52"
$ws.Range("E50").Value2 = "  -2.68%  "
$ws.Range("E51").Value2 = "  +0.08%  "
